$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-04 Thursday" "2025-09-05 Friday"

Replace-Text "810×3=" "623×3="
Replace-Text "649×4=" "127×3="
Replace-Text "278×6=" "924×5="
Replace-Text "650×5=" "554×4="
Replace-Text "413×6=" "112×8="
Replace-Text "857×7=" "627×7="
Replace-Text "257×2=" "194×7="
Replace-Text "422×9=" "357×5="
Replace-Text "166×3=" "388×3="
Replace-Text "171×3=" "329×5="
Replace-Text "302×7=" "777×8="
Replace-Text "810×6=" "801×7="
Replace-Text "751×4=" "425×5="
Replace-Text "995×7=" "643×7="
Replace-Text "533×8=" "667×6="
Replace-Text "106×5=" "821×2="
Replace-Text "413×8=" "734×7="
Replace-Text "455×5=" "521×3="
Replace-Text "346×5=" "438×4="
Replace-Text "767×3=" "288×9="
Replace-Text "492×7=" "666×7="
Replace-Text "227×4=" "919×5="
Replace-Text "610×6=" "766×5="
Replace-Text "906×4=" "754×5="
Replace-Text "905×3=" "361×3="

Write-Output "Done"
